$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> the real Color value for rows whose Color was
# "TOO_LOW_TESTING" before this edit. These rows get a proper color in
# column B and HasLowTesting = TRUE in the new column C.
$lowTestingColors = @{
    10 = "RED"; 21 = "GREEN"; 22 = "ORANGE"; 31 = "RED"; 32 = "RED";
    35 = "ORANGE"; 37 = "GREEN"; 40 = "GREEN"; 43 = "RED"; 46 = "ORANGE";
    56 = "ORANGE"; 67 = "ORANGE"; 68 = "ORANGE"; 70 = "ORANGE"; 75 = "ORANGE";
    77 = "ORANGE"; 80 = "GREEN"; 83 = "ORANGE"; 86 = "ORANGE"; 89 = "RED";
    99 = "ORANGE"; 103 = "GREEN"; 107 = "ORANGE"; 108 = "GREEN"; 112 = "GREEN";
    113 = "GREEN"; 118 = "ORANGE"; 119 = "GREEN"; 120 = "GREEN"
}

# Insert a new column at C; existing C:F (LastUpdate, NotificationRatePer100000,
# PositiveRate, TestsPer100000) shift right to D:G.
$ws.Columns("C").Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "HasLowTesting"

$lastRow = 121
for ($row = 2; $row -le $lastRow; $row++) {
    if ($lowTestingColors.ContainsKey($row)) {
        $ws.Cells.Item($row, 2).Value = $lowTestingColors[$row]
        $ws.Cells.Item($row, 3).Value = $true
    } else {
        $ws.Cells.Item($row, 3).Value = $false
    }
}

# Row 121 is the sentinel "NO_DATA" fallback row (not a real country); it
# keeps HasLowTesting = TRUE regardless of its Color.
$ws.Cells.Item(121, 3).Value = $true
